$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9571118980407642
$ws.Range("D2").Value = 0.3452673588652468

$ws.Range("C3").Value = 0.9691344718981443
$ws.Range("D3").Value = 0.3393220935692614

$ws.Range("C4").Value = 0.3851038752975169
$ws.Range("D4").Value = 0.7025591562145015

$ws.Range("C5").Value = 0.9894889584238222
$ws.Range("D5").Value = 0.3294138930056869

$ws.Range("C6").Value = -0.08030479467229548
$ws.Range("D6").Value = 0.9364656805152671

$ws.Range("C7").Value = -0.8268988214588664
$ws.Range("D7").Value = 0.414060221268185

$ws.Range("C8").Value = -0.009763160005758185
$ws.Range("D8").Value = 0.9922673120589833

$ws.Range("C9").Value = -0.5712722305884617
$ws.Range("D9").Value = 0.5715723857575004

$ws.Range("C10").Value = 0.06660549274157744
$ws.Range("D10").Value = 0.9472859013716244

$ws.Range("C11").Value = 0.5541439180746444
$ws.Range("D11").Value = 0.5831070744628244
